$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75 (shifts existing rows 75-97 down to 76-98)
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row with the new weekly price entry
$ws.Range("A75").Value = 10
$ws.Range("B75").Value = "Vega Modelo de Temuco"
$ws.Range("C75").Value = "La Araucanía"
$ws.Range("D75").Value = 45016
$ws.Range("E75").Value = 9
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100108
$ws.Range("H75").Value = "Tropicales y subtropicales"
$ws.Range("I75").Value = 100108004
$ws.Range("J75").Value = "Papaya"
$ws.Range("K75").Value = "Cultivar IV Región"
$ws.Range("L75").Value = "Primera"
$ws.Range("M75").Value = 55
$ws.Range("N75").Value = 28000
$ws.Range("O75").Value = 28000
$ws.Range("P75").Value = 28000
$ws.Range("Q75").Value = "$/bandeja 10 kilos"
$ws.Range("R75").Value = "Provincia del Elquí"
$ws.Range("S75").Value = 2800
$ws.Range("T75").Value = 10
